# New PO forecast model
# Updates:
#  - "Weekly Quantity" sheet: append rows 61-62
#  - "Monthly Trend" sheet: append row 23
#  - "PO Forecast" sheet: recompute B2:B60, shift/extend rows 61-70

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Weekly Quantity" - append two new weekly observations
# ---------------------------------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")

$weeklyNewDates = @(45662.99999999999, 45676.99999999999)
$weeklyNewQty   = @(90, 50)

$weeklyStartRow = 61
for ($i = 0; $i -lt $weeklyNewDates.Length; $i++) {
    $row = $weeklyStartRow + $i
    $wsWeekly.Cells.Item($row, 1).Value = $weeklyNewDates[$i]
    $wsWeekly.Cells.Item($row, 2).Value = $weeklyNewQty[$i]
}
# Match the date/time number format used by the rest of column A
$wsWeekly.Range("A61:A62").NumberFormat = $wsWeekly.Range("A60").NumberFormat

# ---------------------------------------------------------------------------
# Sheet 2: "Monthly Trend" - append one new monthly observation
# ---------------------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

$wsMonthly.Cells.Item(23, 1).Value = 45688.99999999999
$wsMonthly.Cells.Item(23, 2).Value = 140
$wsMonthly.Range("A23").NumberFormat = $wsMonthly.Range("A22").NumberFormat

# ---------------------------------------------------------------------------
# Sheet 3: "PO Forecast" - new forecast model values
# ---------------------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("PO Forecast")

# Updated forecast quantities for existing rows 2-60 (dates in column A unchanged)
$forecastValues = @(195,199,206,213,220,224,234,238,252,255,259,266,269,273,280,283,287,290,294,297,301,304,308,318,322,325,332,336,339,346,357,360,364,371,375,378,382,385,389,396,399,403,406,410,417,420,434,438,441,466,469,473,476,487,490,504,508,525,529)

$forecastStartRow = 2
for ($i = 0; $i -lt $forecastValues.Length; $i++) {
    $row = $forecastStartRow + $i
    $wsForecast.Cells.Item($row, 2).Value = $forecastValues[$i]
}

# Rows 61-70: dates shift forward one period and extend with two brand new rows
$forecastTailDates = @(45662.99999999999,45676.99999999999,45683.99999999999,45690.99999999999,45697.99999999999,45704.99999999999,45711.99999999999,45718.99999999999,45725.99999999999,45732.99999999999)
$forecastTailValues = @(540,547,550,554,557,561,564,568,571,575)

$forecastTailStartRow = 61
for ($i = 0; $i -lt $forecastTailDates.Length; $i++) {
    $row = $forecastTailStartRow + $i
    $wsForecast.Cells.Item($row, 1).Value = $forecastTailDates[$i]
    $wsForecast.Cells.Item($row, 2).Value = $forecastTailValues[$i]
}
# Match the date/time number format used by the rest of column A
$wsForecast.Range("A61:A70").NumberFormat = $wsForecast.Range("A60").NumberFormat

Write-Host "Done applying new PO forecast model"
